$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 228.2
$ws.Range("I9").Value = 235.875
$ws.Range("J9").Value = 197.5
$ws.Range("K9").Value = 235.875
$ws.Range("L9").Value = 197.5
$ws.Range("M9").Value = -66.875
$ws.Range("N9").Value = -535.5

$ws.Range("H76").Value = 4554.909
$ws.Range("I76").Value = 4567.222
$ws.Range("K76").Value = 4567.222
$ws.Range("M76").Value = -4252.222

$ws.Range("H79").Value = 4554.909
$ws.Range("I79").Value = 4567.222
$ws.Range("K79").Value = 4567.222
$ws.Range("M79").Value = -3475.222

$ws.Range("H86").Value = 7500
$ws.Range("I86").Value = 7500
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 7500
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -6377
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 7500
$ws.Range("I89").Value = 7500
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 37500
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -31884
$ws.Range("N89").ClearContents()

$ws.Range("H129").Value = 1255
$ws.Range("J129").Value = 1570
$ws.Range("L129").Value = 4710
$ws.Range("N129").Value = -14710

$ws.Range("H137").Value = 2176.1714
$ws.Range("I137").Value = 1422.1482
$ws.Range("K137").Value = 4266.444600000001
$ws.Range("M137").Value = -1716.444600000001

$ws.Range("H138").Value = 3604.077
$ws.Range("J138").Value = 3249.5
$ws.Range("L138").Value = 9748.5
$ws.Range("N138").Value = -20028.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4404.385
$ws.Range("I61").Value = 4106.8887
$ws.Range("J61").Value = 5073.75
$ws.Range("K61").Value = 4106.8887
$ws.Range("L61").Value = 5073.75
$ws.Range("M61").Value = -3894.8887
$ws.Range("N61").Value = -5497.75

$ws.Range("H74").Value = 1329.04
$ws.Range("I74").Value = 1240.2609
$ws.Range("K74").Value = 1240.2609
$ws.Range("M74").Value = -366.2609

$ws.Range("H77").Value = 1329.04
$ws.Range("I77").Value = 1240.2609
$ws.Range("K77").Value = 6201.3045
$ws.Range("M77").Value = -1833.3045

$ws.Range("H97").Value = 783.44446
$ws.Range("I97").Value = 1703
$ws.Range("J97").Value = 323.66666
$ws.Range("K97").Value = 1703
$ws.Range("L97").Value = 323.66666
$ws.Range("M97").Value = -1207
$ws.Range("N97").Value = -1315.66666

$ws.Range("H102").Value = 2730.1428
$ws.Range("I102").Value = 3052.75
$ws.Range("K102").Value = 3052.75
$ws.Range("M102").Value = -1430.75

$ws.Range("H110").Value = 2201.5173
$ws.Range("I110").Value = 1616.8
$ws.Range("K110").Value = 1616.8
$ws.Range("M110").Value = 428.2

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H132").Value = 15156805
$ws.Range("I132").Value = 4836.75
$ws.Range("J132").Value = 55562052
$ws.Range("K132").Value = 14510.25
$ws.Range("L132").Value = 166686156
$ws.Range("M132").Value = -11980.25
$ws.Range("N132").Value = -166691216

$ws.Range("H136").Value = 4404.385
$ws.Range("I136").Value = 4106.8887
$ws.Range("J136").Value = 5073.75
$ws.Range("K136").Value = 12320.6661
$ws.Range("L136").Value = 15221.25
$ws.Range("M136").Value = -9770.666100000002
$ws.Range("N136").Value = -20321.25

$ws.Range("H140").Value = 79642.336
$ws.Range("I140").Value = 75000
$ws.Range("J140").Value = 80570.8
$ws.Range("K140").Value = 75000
$ws.Range("L140").Value = 80570.8
$ws.Range("M140").Value = -69820
$ws.Range("N140").Value = -90930.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2519.4736
$ws.Range("J20").Value = 3850
$ws.Range("L20").Value = 3850
$ws.Range("N20").Value = -4344

$ws.Range("H40").Value = 44979.5
$ws.Range("J40").Value = 44979.5
$ws.Range("L40").Value = 44979.5
$ws.Range("N40").Value = -45509.5

$ws.Range("H94").Value = 448
$ws.Range("I94").Value = 399
$ws.Range("K94").Value = 399
$ws.Range("M94").Value = 52

$ws.Range("H97").Value = 4351.3335
$ws.Range("J97").Value = 5000
$ws.Range("L97").Value = 5000
$ws.Range("N97").Value = -6982

$ws.Range("H105").Value = 2938.1667
$ws.Range("I105").Value = 2725.25
$ws.Range("K105").Value = 2725.25
$ws.Range("M105").Value = -978.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 23270.4
$ws.Range("J43").Value = 23270.4
$ws.Range("L43").Value = 23270.4
$ws.Range("N43").Value = -23638.4

$ws.Range("H95").Value = 28779.8
$ws.Range("J95").Value = 28779.8
$ws.Range("L95").Value = 28779.8
$ws.Range("N95").Value = -34271.8

$ws.Range("H101").Value = 23270.4
$ws.Range("J101").Value = 23270.4
$ws.Range("L101").Value = 23270.4
$ws.Range("N101").Value = -29760.4

$ws.Range("H134").Value = 8336948
$ws.Range("I134").Value = 2821
$ws.Range("K134").Value = 8463
$ws.Range("M134").Value = -5928

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 97168.516
$ws.Range("I4").Value = 105.63636
$ws.Range("J4").Value = 334433.34
$ws.Range("K4").Value = 316.90908
$ws.Range("L4").Value = 1003300.02
$ws.Range("M4").Value = -204.90908
$ws.Range("N4").Value = -1003524.02

$ws.Range("H48").Value = 8698.5
$ws.Range("J48").Value = 14499.5
$ws.Range("L48").Value = 43498.5
$ws.Range("N48").Value = -43998.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6979.05
$ws.Range("J70").Value = 6811.6
$ws.Range("L70").Value = 6811.6
$ws.Range("N70").Value = -7351.6

$ws.Range("H73").Value = 6979.05
$ws.Range("J73").Value = 6811.6
$ws.Range("L73").Value = 6811.6
$ws.Range("N73").Value = -8683.6

$ws.Range("H97").Value = 799.25
$ws.Range("I97").Value = 799.25
$ws.Range("K97").Value = 799.25
$ws.Range("M97").Value = -303.25

$ws.Range("H107").Value = 642.7646999999999
$ws.Range("I107").Value = 249
$ws.Range("K107").Value = 249
$ws.Range("M107").Value = 1671

$ws.Range("H132").Value = 3273.75
$ws.Range("J132").Value = 3216
$ws.Range("L132").Value = 9648
$ws.Range("N132").Value = -14708

$ws.Range("H136").Value = 171921.6
$ws.Range("J136").Value = 171921.6
$ws.Range("L136").Value = 515764.8
$ws.Range("N136").Value = -520864.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 4804.625
$ws.Range("I46").Value = 4247.5
$ws.Range("K46").Value = 4247.5
$ws.Range("M46").Value = -4059.5

$ws.Range("H61").Value = 4280.0713
$ws.Range("I61").Value = 3674.7273
$ws.Range("K61").Value = 3674.7273
$ws.Range("M61").Value = -3472.7273

$ws.Range("H68").Value = 3008.1538
$ws.Range("I68").Value = 2882.4546
$ws.Range("J68").Value = 3699.5
$ws.Range("K68").Value = 2882.4546
$ws.Range("L68").Value = 3699.5
$ws.Range("M68").Value = -2133.4546
$ws.Range("N68").Value = -5197.5

$ws.Range("H71").Value = 3008.1538
$ws.Range("I71").Value = 2882.4546
$ws.Range("J71").Value = 3699.5
$ws.Range("K71").Value = 14412.273
$ws.Range("L71").Value = 18497.5
$ws.Range("M71").Value = -10668.273
$ws.Range("N71").Value = -25985.5

$ws.Range("H100").Value = 12000
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()

$ws.Range("H113").Value = 4280.0713
$ws.Range("I113").Value = 3674.7273
$ws.Range("K113").Value = 3674.7273
$ws.Range("M113").Value = -1504.7273

$ws.Range("H136").Value = 125004190
$ws.Range("I136").Value = 3551
$ws.Range("K136").Value = 10653
$ws.Range("M136").Value = -8103

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 15841.667
$ws.Range("I55").Value = 9000
$ws.Range("K55").Value = 9000
$ws.Range("M55").Value = -8723

$ws.Range("H96").Value = 3800.375
$ws.Range("I96").Value = 9501.5
$ws.Range("K96").Value = 9501.5
$ws.Range("M96").Value = -8128.5

$ws.Range("H100").Value = 1686.8334
$ws.Range("I100").Value = 1061.3636
$ws.Range("K100").Value = 2122.7272
$ws.Range("M100").Value = -1581.7272
